$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 152
$ws.Range("I11").Value = 152
$ws.Range("K11").Value = 152
$ws.Range("M11").Value = -12
$ws.Range("H15").Value = 1964.6615
$ws.Range("I15").Value = 1964.6615
$ws.Range("K15").Value = 5893.9845
$ws.Range("M15").Value = -5724.9845
$ws.Range("H38").Value = 5977.905
$ws.Range("I38").Value = 1129.8182
$ws.Range("K38").Value = 3389.4546
$ws.Range("M38").Value = -3017.4546
$ws.Range("H41").Value = 457
$ws.Range("I41").Value = 200
$ws.Range("K41").Value = 200
$ws.Range("M41").Value = 240
$ws.Range("H47").Value = 200
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = $null
$ws.Range("H54").Value = 9969.333000000001
$ws.Range("I54").Value = 9954
$ws.Range("J54").Value = 9977
$ws.Range("K54").Value = 9954
$ws.Range("L54").Value = 9977
$ws.Range("M54").Value = -9468
$ws.Range("N54").Value = -10949
$ws.Range("H86").Value = 2416.5557
$ws.Range("I86").Value = 3454
$ws.Range("J86").Value = 1586.6
$ws.Range("K86").Value = 3454
$ws.Range("L86").Value = 1586.6
$ws.Range("M86").Value = -2331
$ws.Range("N86").Value = -3832.6
$ws.Range("H89").Value = 2416.5557
$ws.Range("I89").Value = 3454
$ws.Range("J89").Value = 1586.6
$ws.Range("K89").Value = 17270
$ws.Range("L89").Value = 7933
$ws.Range("M89").Value = -11654
$ws.Range("N89").Value = -19165
$ws.Range("H100").Value = 3097.5908
$ws.Range("I100").Value = 1977.1333
$ws.Range("K100").Value = 1977.1333
$ws.Range("M100").Value = -1436.1333
$ws.Range("H106").Value = 48469.285
$ws.Range("I106").Value = 58856.9
$ws.Range("K106").Value = 58856.9
$ws.Range("M106").Value = -58225.9
$ws.Range("H107").Value = 4000
$ws.Range("I107").Value = 4000
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 4000
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = -2080
$ws.Range("N107").Value = -7840
$ws.Range("H118").Value = 1392.2106
$ws.Range("I118").Value = 914
$ws.Range("K118").Value = 2742
$ws.Range("M118").Value = -1085
$ws.Range("H132").Value = 3042.361
$ws.Range("I132").Value = 2703.875
$ws.Range("J132").Value = 5750.25
$ws.Range("K132").Value = 8111.625
$ws.Range("L132").Value = 17250.75
$ws.Range("M132").Value = -5581.625
$ws.Range("N132").Value = -22310.75
$ws.Range("H135").Value = 490.31708
$ws.Range("I135").Value = 495.58975
$ws.Range("K135").Value = 4460.30775
$ws.Range("M135").Value = -1925.30775
$ws.Range("H137").Value = 26322842
$ws.Range("I137").Value = 38470744
$ws.Range("J137").Value = 2385.8333
$ws.Range("K137").Value = 115412232
$ws.Range("L137").Value = 7157.499899999999
$ws.Range("M137").Value = -115409682
$ws.Range("N137").Value = -12257.4999
$ws.Range("H138").Value = 3472.2424
$ws.Range("I138").Value = 1631.5946
$ws.Range("J138").Value = 4570.6934
$ws.Range("K138").Value = 4894.783799999999
$ws.Range("L138").Value = 13712.0802
$ws.Range("M138").Value = 245.2162000000008
$ws.Range("N138").Value = -23992.0802
$ws.Range("H141").Value = 6390.7046
$ws.Range("I141").Value = 3491.7222
$ws.Range("K141").Value = 10475.1666
$ws.Range("M141").Value = -5295.1666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 2750
$ws.Range("I38").Value = 2750
$ws.Range("K38").Value = 2750
$ws.Range("M38").Value = -2283
$ws.Range("H45").Value = 1386
$ws.Range("I45").Value = 1386
$ws.Range("K45").Value = 1386
$ws.Range("M45").Value = -1009
$ws.Range("H61").Value = 1256.3541
$ws.Range("I61").Value = 1200.1277
$ws.Range("K61").Value = 1200.1277
$ws.Range("M61").Value = -988.1277
$ws.Range("H63").Value = 1427.4286
$ws.Range("I63").Value = 1548.8334
$ws.Range("K63").Value = 1548.8334
$ws.Range("M63").Value = -862.8334
$ws.Range("H66").Value = 1427.4286
$ws.Range("I66").Value = 1548.8334
$ws.Range("K66").Value = 7744.166999999999
$ws.Range("M66").Value = -4312.166999999999
$ws.Range("H74").Value = 10826.366
$ws.Range("I74").Value = 7183.3335
$ws.Range("K74").Value = 7183.3335
$ws.Range("M74").Value = -6309.3335
$ws.Range("H77").Value = 10826.366
$ws.Range("I77").Value = 7183.3335
$ws.Range("K77").Value = 35916.6675
$ws.Range("M77").Value = -31548.6675
$ws.Range("H94").Value = 69711.06
$ws.Range("J94").Value = 70598.39999999999
$ws.Range("L94").Value = 70598.39999999999
$ws.Range("N94").Value = -72400.39999999999
$ws.Range("H95").Value = 29966.666
$ws.Range("J95").Value = 29966.666
$ws.Range("L95").Value = 29966.666
$ws.Range("N95").Value = -35458.666
$ws.Range("H97").Value = 511.16
$ws.Range("I97").Value = 525.55554
$ws.Range("K97").Value = 525.55554
$ws.Range("M97").Value = -29.55553999999995
$ws.Range("H122").Value = 4534.8857
$ws.Range("I122").Value = 4621.206
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 13863.618
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -11413.618
$ws.Range("N122").Value = -9700
$ws.Range("H132").Value = 9171.777
$ws.Range("I132").Value = 6208.4443
$ws.Range("K132").Value = 18625.3329
$ws.Range("M132").Value = -16095.3329
$ws.Range("H136").Value = 1256.3541
$ws.Range("I136").Value = 1200.1277
$ws.Range("K136").Value = 3600.3831
$ws.Range("M136").Value = -1050.3831

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 32260908
$ws.Range("I20").Value = 50002156
$ws.Range("K20").Value = 50002156
$ws.Range("M20").Value = -50001909
$ws.Range("H42").Value = 119980
$ws.Range("J42").Value = 119980
$ws.Range("L42").Value = 119980
$ws.Range("N42").Value = -120636
$ws.Range("H86").Value = 2870.9148
$ws.Range("I86").Value = 2605.9443
$ws.Range("K86").Value = 2605.9443
$ws.Range("M86").Value = -1482.9443
$ws.Range("H89").Value = 2870.9148
$ws.Range("I89").Value = 2605.9443
$ws.Range("K89").Value = 13029.7215
$ws.Range("M89").Value = -7413.7215
$ws.Range("H92").Value = 81289.8
$ws.Range("J92").Value = 81289.8
$ws.Range("L92").Value = 81289.8
$ws.Range("N92").Value = -86281.8
$ws.Range("H94").Value = 2327.1843
$ws.Range("I94").Value = 1280.8
$ws.Range("J94").Value = 4339.4614
$ws.Range("K94").Value = 1280.8
$ws.Range("L94").Value = 4339.4614
$ws.Range("M94").Value = -829.8
$ws.Range("N94").Value = -5241.4614
$ws.Range("H105").Value = 2306.3333
$ws.Range("I105").Value = 2178.6155
$ws.Range("K105").Value = 2178.6155
$ws.Range("M105").Value = -431.6154999999999
$ws.Range("H107").Value = 3201.6875
$ws.Range("I107").Value = 3201.6875
$ws.Range("K107").Value = 3201.6875
$ws.Range("M107").Value = -1281.6875
$ws.Range("H134").Value = 1554.0597
$ws.Range("I134").Value = 1564.9517
$ws.Range("J134").Value = 1419
$ws.Range("K134").Value = 4694.855100000001
$ws.Range("L134").Value = 4257
$ws.Range("M134").Value = -2159.855100000001
$ws.Range("N134").Value = -9327

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1850.6364
$ws.Range("I12").Value = 1222.4286
$ws.Range("K12").Value = 1222.4286
$ws.Range("M12").Value = -1052.4286
$ws.Range("H31").Value = 241551.81
$ws.Range("I31").Value = 4369.3335
$ws.Range("J31").Value = 1003924.06
$ws.Range("K31").Value = 4369.3335
$ws.Range("L31").Value = 1003924.06
$ws.Range("M31").Value = -4074.3335
$ws.Range("N31").Value = -1004514.06
$ws.Range("H34").Value = 241551.81
$ws.Range("I34").Value = 4369.3335
$ws.Range("J34").Value = 1003924.06
$ws.Range("K34").Value = 4369.3335
$ws.Range("L34").Value = 1003924.06
$ws.Range("M34").Value = -4167.3335
$ws.Range("N34").Value = -1004328.06
$ws.Range("H58").Value = 2592.2727
$ws.Range("I58").Value = 2913.1
$ws.Range("K58").Value = 2913.1
$ws.Range("M58").Value = -2710.1
$ws.Range("H94").Value = 958.8946999999999
$ws.Range("I94").Value = 918.5714
$ws.Range("K94").Value = 918.5714
$ws.Range("M94").Value = -467.5714
$ws.Range("H107").Value = 1160.4783
$ws.Range("I107").Value = 558.1667
$ws.Range("K107").Value = 558.1667
$ws.Range("M107").Value = 1361.8333
$ws.Range("H132").Value = 2310.9033
$ws.Range("I132").Value = 2251.7727
$ws.Range("K132").Value = 6755.3181
$ws.Range("M132").Value = -4225.3181
$ws.Range("H136").Value = 2592.2727
$ws.Range("I136").Value = 2913.1
$ws.Range("K136").Value = 8739.299999999999
$ws.Range("M136").Value = -6189.299999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 10101204
$ws.Range("I2").Value = 115.6
$ws.Range("J2").Value = 18518778
$ws.Range("K2").Value = 693.5999999999999
$ws.Range("L2").Value = 111112668
$ws.Range("M2").Value = -580.5999999999999
$ws.Range("N2").Value = -111112894
$ws.Range("H5").Value = 758.0333000000001
$ws.Range("I5").Value = 606.7727
$ws.Range("J5").Value = 1174
$ws.Range("K5").Value = 1820.3181
$ws.Range("L5").Value = 3522
$ws.Range("M5").Value = -1708.3181
$ws.Range("N5").Value = -3746
$ws.Range("H33").Value = 623.125
$ws.Range("I33").Value = 521.25
$ws.Range("J33").Value = 725
$ws.Range("K33").Value = 3127.5
$ws.Range("L33").Value = 4350
$ws.Range("M33").Value = -2844.5
$ws.Range("N33").Value = -4916
$ws.Range("H48").Value = 9999
$ws.Range("J48").Value = 9999
$ws.Range("L48").Value = 29997
$ws.Range("N48").Value = -30497
$ws.Range("H55").Value = 3997.6
$ws.Range("J55").Value = 3663
$ws.Range("L55").Value = 10989
$ws.Range("N55").Value = -11343
$ws.Range("H64").Value = 2249.5
$ws.Range("I64").Value = 2249.5
$ws.Range("K64").Value = 6748.5
$ws.Range("M64").Value = -6478.5
$ws.Range("H67").Value = 2249.5
$ws.Range("I67").Value = 2249.5
$ws.Range("K67").Value = 6748.5
$ws.Range("M67").Value = -5812.5
$ws.Range("H122").Value = 1108.0952
$ws.Range("I122").Value = 2290.5
$ws.Range("J122").Value = 635.13336
$ws.Range("K122").Value = 20614.5
$ws.Range("L122").Value = 5716.20024
$ws.Range("M122").Value = -18164.5
$ws.Range("N122").Value = -10616.20024
$ws.Range("H135").Value = 758.0333000000001
$ws.Range("I135").Value = 606.7727
$ws.Range("J135").Value = 1174
$ws.Range("K135").Value = 5460.954299999999
$ws.Range("L135").Value = 10566
$ws.Range("M135").Value = -2925.954299999999
$ws.Range("N135").Value = -15636
$ws.Range("H137").Value = 2765.4
$ws.Range("I137").Value = 1539.2
$ws.Range("J137").Value = 3991.6
$ws.Range("K137").Value = 4617.6
$ws.Range("L137").Value = 11974.8
$ws.Range("M137").Value = 482.3999999999996
$ws.Range("N137").Value = -22174.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
$ws.Range("H28").Value = 19282.363
$ws.Range("J28").Value = 19282.363
$ws.Range("L28").Value = 19282.363
$ws.Range("N28").Value = -19666.363
$ws.Range("H55").Value = 21748.75
$ws.Range("I55").Value = 3000
$ws.Range("K55").Value = 3000
$ws.Range("M55").Value = -2673
$ws.Range("H70").Value = 22382.34
$ws.Range("I70").Value = 33772.668
$ws.Range("K70").Value = 33772.668
$ws.Range("M70").Value = -33502.668
$ws.Range("H73").Value = 22382.34
$ws.Range("I73").Value = 33772.668
$ws.Range("K73").Value = 33772.668
$ws.Range("M73").Value = -32836.668
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
$ws.Range("H97").Value = 1057.1428
$ws.Range("I97").Value = 1060.5454
$ws.Range("K97").Value = 1060.5454
$ws.Range("M97").Value = -564.5454
$ws.Range("H102").Value = 1454.123
$ws.Range("I102").Value = 1388.74
$ws.Range("K102").Value = 1388.74
$ws.Range("M102").Value = 233.26
$ws.Range("H132").Value = 4503.965
$ws.Range("I132").Value = 3485.6597
$ws.Range("J132").Value = 9290
$ws.Range("K132").Value = 10456.9791
$ws.Range("L132").Value = 27870
$ws.Range("M132").Value = -7926.9791
$ws.Range("N132").Value = -32930

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 58829320
$ws.Range("I7").Value = 111113170
$ws.Range("K7").Value = 111113170
$ws.Range("M7").Value = -111113058
$ws.Range("H40").Value = 3350.6875
$ws.Range("I40").Value = 3279.1428
$ws.Range("J40").Value = 3851.5
$ws.Range("K40").Value = 3279.1428
$ws.Range("L40").Value = 3851.5
$ws.Range("M40").Value = -3143.1428
$ws.Range("N40").Value = -4123.5
$ws.Range("H45").Value = 20873.875
$ws.Range("I45").Value = 17832.166
$ws.Range("K45").Value = 17832.166
$ws.Range("M45").Value = -17425.166
$ws.Range("H46").Value = 689.75
$ws.Range("J46").Value = 690.3333
$ws.Range("L46").Value = 690.3333
$ws.Range("N46").Value = -1066.3333
$ws.Range("H97").Value = 78330
$ws.Range("J97").Value = 78330
$ws.Range("L97").Value = 78330
$ws.Range("N97").Value = -80312
$ws.Range("H126").Value = 58829320
$ws.Range("I126").Value = 111113170
$ws.Range("K126").Value = 333339510
$ws.Range("M126").Value = -333337040
$ws.Range("H132").Value = 2622.5898
$ws.Range("I132").Value = 2428.4473
$ws.Range("K132").Value = 7285.341899999999
$ws.Range("M132").Value = -4755.341899999999
$ws.Range("H133").Value = 62249.75
$ws.Range("J133").Value = 62249.75
$ws.Range("L133").Value = 62249.75
$ws.Range("N133").Value = -67309.75
$ws.Range("H136").Value = 5165.407
$ws.Range("I136").Value = 4938.64
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 14815.92
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -12265.92
$ws.Range("N136").Value = -29100

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1133.2
$ws.Range("I107").Value = 992.125
$ws.Range("J107").Value = 1697.5
$ws.Range("K107").Value = 2976.375
$ws.Range("L107").Value = 5092.5
$ws.Range("M107").Value = -1056.375
$ws.Range("N107").Value = -8932.5
$ws.Range("H126").Value = 14494081
$ws.Range("I126").Value = 19609224
$ws.Range("J126").Value = 1172.3334
$ws.Range("K126").Value = 58827672
$ws.Range("L126").Value = 3517.0002
$ws.Range("M126").Value = -58825202
$ws.Range("N126").Value = -8457.0002
$ws.Range("H132").Value = 1544.2307
$ws.Range("I132").Value = 1314.7241
$ws.Range("J132").Value = 2209.8
$ws.Range("K132").Value = 3944.1723
$ws.Range("L132").Value = 6629.400000000001
$ws.Range("M132").Value = -1414.1723
$ws.Range("N132").Value = -11689.4
$ws.Range("H136").Value = 5978.2085
$ws.Range("J136").Value = 8999
$ws.Range("L136").Value = 26997
$ws.Range("N136").Value = -32097
